# Remove old, missing files (LICENSE, README.md, index.html) from the
# "hasPart" entry of the RootDataset entity. The regenerating tool drops
# the stale "hasPart" property and re-appends the recomputed one at the
# end of the property list, so mirror that by deleting the existing
# "hasPart" row and writing a fresh one after the last row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RootDataset")

# Locate the "hasPart" row dynamically (column A holds property names,
# column B holds the corresponding values) rather than assuming a fixed
# row number.
$hasPartRow = -1
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value()
    if ($name -eq "hasPart") {
        $hasPartRow = $r
        break
    }
}

if ($hasPartRow -eq -1) {
    throw "Could not locate 'hasPart' row on RootDataset sheet"
}

$ws.Rows.Item($hasPartRow).Delete()

$newRow = $lastRow
$ws.Cells.Item($newRow, 1).Value = "hasPart"
$ws.Cells.Item($newRow, 2).Value = '["lots_of_little_files/", "pics/"]'
